$wb = $excel.ActiveWorkbook

# Add a new worksheet that documents the workbook and each data tab, and
# make sure it ends up as the very first sheet (tab order becomes:
# Documentation, AVA, AF_Trans, FILL Table).
$docSheet = $wb.Worksheets.Add()
$docSheet.Name = "Documentation"
$docSheet.Move($wb.Worksheets.Item(1))

$docSheet.Range("A1").Value = "Workbook: Mapping and transformation of new process in Industrial sector across regions"
$docSheet.Range("A2").Value = "AVA: Process availability across regions"
$docSheet.Range("A3").Value = "AF_Trans: Time-slice-specific availability factors for existing energy processes"
$docSheet.Range("A4").Value = "FILL Table: Model data for transformation operations"
